$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F, applied identically to both
# the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.
$updates = @{
    3  = 24
    4  = 969
    6  = 2326
    8  = 1425
    12 = 379
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
